{"js": "// Add the transaction-date placeholder to the invoice item row.\n//\n// The template has two copies of the same table (paginated invoice);\n// each copy has a cell with the tbs placeholder\n//   [detailInsurance_sub1.transactionDate; ope=formatdate; format='DD-mm-YYYY'; block=tbs:row;]\n// Both copies need \"mm\" capitalized to \"MM\". In the first table the\n// placeholder is split across several runs (the cell text was typed in\n// pieces) and needs to collapse into a single run. In the second table\n// the placeholder is already a single run but a stray \"_GoBack\" bookmark\n// (left over from the last cursor position when the file was saved in\n// Word) needs to move from the very end of the document into the middle\n// of that run, splitting it in two.\n\nconst body = context.document.body;\n\nconst oldText =\n  \"[detailInsurance_sub1.transactionDate; ope=formatdate; format=\\u2019DD-mm-YYYY\\u2019; block=tbs:row;]\";\nconst newText =\n  \"[detailInsurance_sub1.transactionDate; ope=formatdate; format=\\u2019DD-MM-YYYY\\u2019; block=tbs:row;]\";\n\n// 1. Drop the old \"_GoBack\" bookmark sitting in the trailing empty\n//    paragraph at the end of the document (it gets reinserted below, in\n//    the middle of the second table's placeholder text).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Normalize the placeholder text in both tables (this also merges the\n//    first table's multi-run cell into a single run, since Range.insertText\n//    with \"Replace\" rewrites the whole matched span as one run).\nconst matches = body.search(oldText, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < matches.items.length; i++) {\n  matches.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n\n// 3. Re-insert the \"_GoBack\" bookmark in the SECOND table's placeholder,\n//    right after \"format='DD-MM\" (splitting that run into two runs with\n//    the bookmark sandwiched between them, matching the original cursor\n//    position Word recorded there).\nconst splitMatches = body.search(\"format=\\u2019DD-MM\", { matchCase: true });\nsplitMatches.load(\"items\");\nawait context.sync();\n\nconst splitPoint = splitMatches.items[splitMatches.items.length - 1].getRange(\"After\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Add the transaction-date placeholder to the invoice item row.\n#\n# The template has two copies of the same table (paginated invoice);\n# each copy has a cell with the tbs placeholder\n#   [detailInsurance_sub1.transactionDate; ope=formatdate; format='DD-mm-YYYY'; block=tbs:row;]\n# Both copies need \"mm\" capitalized to \"MM\". In the first table the\n# placeholder is split across several runs (the cell text was typed in\n# pieces) and needs to collapse into a single run. In the second table\n# the placeholder is already a single run but a stray \"_GoBack\" bookmark\n# (left over from the last cursor position when the file was saved in\n# Word) needs to move from the very end of the document into the middle\n# of that run, splitting it in two.\n\n$d = $word.ActiveDocument\n$quote = [char]8217\n\n# 1. Drop the old \"_GoBack\" bookmark sitting in the trailing empty\n#    paragraph at the end of the document (it gets reinserted below, in\n#    the middle of the second table's placeholder text).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Normalize the placeholder text in both tables (Find/Replace across\n#    the whole story rewrites each matched span as a single run, which\n#    also merges the first table's multi-run cell into one run).\n$oldText = \"[detailInsurance_sub1.transactionDate; ope=formatdate; format=\" + $quote + \"DD-mm-YYYY\" + $quote + \"; block=tbs:row;]\"\n$newText = \"[detailInsurance_sub1.transactionDate; ope=formatdate; format=\" + $quote + \"DD-MM-YYYY\" + $quote + \"; block=tbs:row;]\"\n\n$replaceRange = $d.Content\n$find = $replaceRange.Find\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.MatchCase = $true\n$find.Execute([ref]$find.Text, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2) | Out-Null\n\n# 3. Re-insert the \"_GoBack\" bookmark in the SECOND table's placeholder,\n#    right after \"format='DD-MM\" (splitting that run into two runs with\n#    the bookmark sandwiched between them, matching the original cursor\n#    position Word recorded there). Walk every match so we land on the\n#    last (second) occurrence.\n$splitMarker = \"format=\" + $quote + \"DD-MM\"\n$searchRange = $d.Content\n$searchFind = $searchRange.Find\n$searchFind.Text = $splitMarker\n$searchFind.MatchCase = $true\n$searchFind.Forward = $true\n$splitPos = -1\nwhile ($searchFind.Execute()) {\n    $splitPos = $searchRange.End\n    $searchRange.Collapse(0)\n}\n\n$bmRange = $d.Range($splitPos, $splitPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
